$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new pollutant rows: SO2 (row 6) and CO (row 7)
$ws.Range("A6").Value = "SO2"
$ws.Range("B6").Value = 30
$ws.Range("C6").Value = 100
$ws.Range("D6").Value = "NA"

$ws.Range("A7").Value = "CO"
$ws.Range("B7").Value = "NA"
$ws.Range("C7").Value = 8
$ws.Range("D7").Value = "NA"

# Match styling of existing data rows (right-aligned values, like B2:D5)
# Row 6: only B6 and D6 are right-aligned (C6 keeps default alignment)
$ws.Range("B6").HorizontalAlignment = -4152
$ws.Range("D6").HorizontalAlignment = -4152
$ws.Range("B7:D7").HorizontalAlignment = -4152

# Update selection to match the new active cell recorded in the diff
$ws.Range("H13").Select()
